# Generate Report for Handback
# Refresh the handoff/handback timestamps (and priority "ht" -> "mt")
# for the ae14b882-... and d577b113-... entries across the Overview,
# zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("G4").Value = "2016-08-13 10:19:15"
$ws.Range("G5").Value = "2016-08-13 10:19:15"

# --- zh-cn sheet --------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("E4").Value = "mt"
$ws.Range("E5").Value = "mt"
$ws.Range("H4").Value = "2016-08-13 10:19:08"
$ws.Range("H5").Value = "2016-08-13 10:19:08"
$ws.Range("K4").Value = "2016-08-13 10:19:36"
$ws.Range("K5").Value = "2016-08-13 10:19:36"

# --- de-de sheet --------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("E4").Value = "mt"
$ws.Range("E5").Value = "mt"
$ws.Range("H4").Value = "2016-08-13 10:19:15"
$ws.Range("H5").Value = "2016-08-13 10:19:15"
$ws.Range("K4").Value = "2016-08-13 10:19:46"
$ws.Range("K5").Value = "2016-08-13 10:19:46"
